$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.606.51"
$ws.Range("E2").Value = "  +3.55%  "
$ws.Range("D3").Value = "2.433.97"
$ws.Range("E3").Value = "  +2.36%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.46%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.36"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.62%  "
$ws.Range("E7").Value = "  +1.27%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  +1.43%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.46"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.59%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0800"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.47%  "
$ws.Range("E12").Value = "  +0.77%  "
$ws.Range("E13").Value = "  +2.22%  "
$ws.Range("E14").Value = "  +2.44%  "
$ws.Range("D15").Value = "2.816.67"
$ws.Range("E15").Value = "  +2.42%  "
$ws.Range("D16").Value = "2.430.78"
$ws.Range("E16").Value = "  +1.33%  "
$ws.Range("E17").Value = "  +4.21%  "
$ws.Range("D18").Value = "44.529.90"
$ws.Range("E18").Value = "  +3.41%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.42"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.70%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.40"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.41%  "
$ws.Range("D21").Value = "0.0₃0908"
$ws.Range("E21").Value = "  +2.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.84"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.96%  "
$ws.Range("E23").Value = "  +4.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "240.78"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.49"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.39%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.24"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.07%  "
$ws.Range("E28").Value = "  -2.54%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.69"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.90%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "33.40"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.14%  "
$ws.Range("B31").Value = "Kaspa"
$ws.Range("C31").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.123"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +16.63%  "
$ws.Range("B32").Value = "Celestia"
$ws.Range("C32").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.54"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +11.89%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.19"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.30%  "
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0763"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.73%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.91"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.15%  "
$ws.Range("E37").Value = "  +3.80%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.90"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.95%  "
$ws.Range("B39").Value = "WEMIXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.34"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.23%  "
$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "126.80"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +11.28%  "
$ws.Range("E41").Value = "  +0.69%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.94"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.25%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0291"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.57%  "
$ws.Range("D44").Value = "1.948.70"
$ws.Range("E44").Value = "  -0.19%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.17"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.15%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.94"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.59"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.80%  "
$ws.Range("E48").Value = "  +10.35%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "53.47"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.95%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "74.02"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.00%  "
$ws.Range("E51").Value = "  +4.94%  "
